$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete data rows (rows 3 and 4), keeping only the header
# row and a single data row (row 2).
$ws.Rows("3:4").Delete()

# Update the remaining data row with the new official title record.
$ws.Range("A2").Value = 803728
$ws.Range("B2").Value = "宣慰使司都元帥府參謀"

# c_dy ("18") and c_source ("2067") must stay text, not be auto-converted to
# numbers, so force a text format before assignment and then clear the
# formatting change back off so no extra cell style is left behind.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "18"
$ws.Range("C2").ClearFormats()

$ws.Range("D2").Value = "Counselor of the Chief Military Command"
$ws.Range("E2").Value = "xuan wei shi si dou yuan shuai fu can mou"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "2067"
$ws.Range("F2").ClearFormats()
